# Weekly update: insert a new "Poroto verde" price record (Feria Lagunitas
# de Puerto Montt) dated 2023-07-27, pushing the existing historical rows
# (old rows 116-160) down by one to new rows 117-161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116; Excel shifts rows 116:160 -> 117:161.
$ws.Rows(116).Insert()

# Populate the newly inserted row 116 with this week's record.
$ws.Cells.Item(116, 1).Value = 4
$ws.Cells.Item(116, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(116, 3).Value = "Los Lagos"
$ws.Cells.Item(116, 4).Value = 45134
$ws.Cells.Item(116, 5).Value = 10
$ws.Cells.Item(116, 6).Value = 100112031
$ws.Cells.Item(116, 7).Value = "Poroto verde"
$ws.Cells.Item(116, 8).Value = "Magnum"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 40
$ws.Cells.Item(116, 11).Value = 29000
$ws.Cells.Item(116, 12).Value = 29000
$ws.Cells.Item(116, 13).Value = 29000
$ws.Cells.Item(116, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(116, 15).Value = "Perú"
$ws.Cells.Item(116, 16).Value = 1160
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"
